$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.986.59'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.274.31'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.88'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.636'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.50'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +2.59%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.449'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +7.09%  '
$ws.Range("E10").Value = '  +7.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.56'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.29'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +14.60%  '
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("D14").Value = '2.613.52'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.11'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +6.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.838'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("D18").Value = '2.258.30'
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").Value = '43.899.91'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("E20").Value = '  +7.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.74'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("E22").Value = '  -2.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.24'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  -4.39%  '
$ws.Range("E26").Value = '  -4.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.06'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.30'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +24.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.60'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +0.98%  '
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.92'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("E32").Value = '  -6.08%  '
$ws.Range("E33").Value = '  +2.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0703'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +6.58%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.80'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +3.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.52'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  -5.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0259'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000228'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +4.96%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0988'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +1.56%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.60'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +5.45%  '
$ws.Range("E45").Value = '  -6.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.43'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +12.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.25'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("E49").Value = '  -4.21%  '
$ws.Range("D50").Value = '1.446.83'
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("E51").Value = '  +0.63%  '
